# Generated Excel COM-interop script applying the Bismarck_Profits diff
$wb = $excel.ActiveWorkbook

### Sheet: ALC ###
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 585.7143   # H5
$ws.Cells.Item(5, 9).Value = 640   # I5
$ws.Cells.Item(5, 10).Value = 450   # J5
$ws.Cells.Item(5, 11).Value = 640   # K5
$ws.Cells.Item(5, 12).Value = 450   # L5
$ws.Cells.Item(5, 13).Value = -525   # M5
$ws.Cells.Item(5, 14).Value = -680   # N5
$ws.Cells.Item(19, 8).Value = 2527.6667   # H19
$ws.Cells.Item(19, 9).Value = 2585.3333   # I19
$ws.Cells.Item(19, 10).Value = 2498.8333   # J19
$ws.Cells.Item(19, 11).Value = 2585.3333   # K19
$ws.Cells.Item(19, 12).Value = 2498.8333   # L19
$ws.Cells.Item(19, 13).Value = -2410.3333   # M19
$ws.Cells.Item(19, 14).Value = -2848.8333   # N19
$ws.Cells.Item(116, 8).Value = 2006   # H116
$ws.Cells.Item(116, 9).Value = 0   # I116
$ws.Cells.Item(116, 10).Value = 2006   # J116
$ws.Cells.Item(116, 11).Value = 0   # K116
$ws.Cells.Item(116, 12).Value = 2006   # L116
$ws.Cells.Item(116, 14).Value = -8890   # N116

### Sheet: CRP ###
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(25, 8).Value = 1250   # H25
$ws.Cells.Item(25, 9).Value = 500   # I25
$ws.Cells.Item(25, 10).Value = 2000   # J25
$ws.Cells.Item(25, 11).Value = 500   # K25
$ws.Cells.Item(25, 12).Value = 2000   # L25
$ws.Cells.Item(25, 13).Value = -326   # M25
$ws.Cells.Item(25, 14).Value = -2348   # N25

### Sheet: GSM ###
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(125, 8).Value = 0   # H125
$ws.Cells.Item(125, 9).Value = 0   # I125
$ws.Cells.Item(125, 10).Value = 0   # J125
$ws.Cells.Item(125, 11).Value = 0   # K125
$ws.Cells.Item(125, 12).Value = 0   # L125
$ws.Cells.Item(126, 8).Value = 0   # H126
$ws.Cells.Item(126, 9).Value = 0   # I126
$ws.Cells.Item(126, 10).Value = 0   # J126
$ws.Cells.Item(126, 11).Value = 0   # K126
$ws.Cells.Item(126, 12).Value = 0   # L126
$ws.Cells.Item(127, 8).Value = 0   # H127
$ws.Cells.Item(127, 9).Value = 0   # I127
$ws.Cells.Item(127, 10).Value = 0   # J127
$ws.Cells.Item(127, 11).Value = 0   # K127
$ws.Cells.Item(127, 12).Value = 0   # L127
$ws.Cells.Item(128, 8).Value = 0   # H128
$ws.Cells.Item(128, 9).Value = 0   # I128
$ws.Cells.Item(128, 10).Value = 0   # J128
$ws.Cells.Item(128, 11).Value = 0   # K128
$ws.Cells.Item(128, 12).Value = 0   # L128
$ws.Cells.Item(129, 8).Value = 0   # H129
$ws.Cells.Item(129, 9).Value = 0   # I129
$ws.Cells.Item(129, 10).Value = 0   # J129
$ws.Cells.Item(129, 11).Value = 0   # K129
$ws.Cells.Item(129, 12).Value = 0   # L129
$ws.Cells.Item(130, 8).Value = 0   # H130
$ws.Cells.Item(130, 9).Value = 0   # I130
$ws.Cells.Item(130, 10).Value = 0   # J130
$ws.Cells.Item(130, 11).Value = 0   # K130
$ws.Cells.Item(130, 12).Value = 0   # L130
$ws.Cells.Item(131, 8).Value = 0   # H131
$ws.Cells.Item(131, 9).Value = 0   # I131
$ws.Cells.Item(131, 10).Value = 0   # J131
$ws.Cells.Item(131, 11).Value = 0   # K131
$ws.Cells.Item(131, 12).Value = 0   # L131
$ws.Cells.Item(132, 8).Value = 3076.7693   # H132
$ws.Cells.Item(132, 9).Value = 2166.6667   # I132
$ws.Cells.Item(132, 10).Value = 3856.8572   # J132
$ws.Cells.Item(132, 11).Value = 6500.000100000001   # K132
$ws.Cells.Item(132, 12).Value = 11570.5716   # L132
$ws.Cells.Item(132, 13).Value = -3970.000100000001   # M132
$ws.Cells.Item(132, 14).Value = -16630.5716   # N132
$ws.Cells.Item(133, 8).Value = 0   # H133
$ws.Cells.Item(133, 9).Value = 0   # I133
$ws.Cells.Item(133, 10).Value = 0   # J133
$ws.Cells.Item(133, 11).Value = 0   # K133
$ws.Cells.Item(133, 12).Value = 0   # L133
$ws.Cells.Item(134, 8).Value = 0   # H134
$ws.Cells.Item(134, 9).Value = 0   # I134
$ws.Cells.Item(134, 10).Value = 0   # J134
$ws.Cells.Item(134, 11).Value = 0   # K134
$ws.Cells.Item(134, 12).Value = 0   # L134
$ws.Cells.Item(135, 8).Value = 0   # H135
$ws.Cells.Item(135, 9).Value = 0   # I135
$ws.Cells.Item(135, 10).Value = 0   # J135
$ws.Cells.Item(135, 11).Value = 0   # K135
$ws.Cells.Item(135, 12).Value = 0   # L135
$ws.Cells.Item(136, 8).Value = 49500   # H136
$ws.Cells.Item(136, 9).Value = 0   # I136
$ws.Cells.Item(136, 10).Value = 49500   # J136
$ws.Cells.Item(136, 11).Value = 0   # K136
$ws.Cells.Item(136, 12).Value = 148500   # L136
$ws.Cells.Item(136, 14).Value = -153600   # N136
$ws.Cells.Item(137, 8).Value = 150000   # H137
$ws.Cells.Item(137, 9).Value = 0   # I137
$ws.Cells.Item(137, 10).Value = 150000   # J137
$ws.Cells.Item(137, 11).Value = 0   # K137
$ws.Cells.Item(137, 12).Value = 150000   # L137
$ws.Cells.Item(137, 14).Value = -160200   # N137
$ws.Cells.Item(138, 8).Value = 0   # H138
$ws.Cells.Item(138, 9).Value = 0   # I138
$ws.Cells.Item(138, 10).Value = 0   # J138
$ws.Cells.Item(138, 11).Value = 0   # K138
$ws.Cells.Item(138, 12).Value = 0   # L138
$ws.Cells.Item(139, 8).Value = 0   # H139
$ws.Cells.Item(139, 9).Value = 0   # I139
$ws.Cells.Item(139, 10).Value = 0   # J139
$ws.Cells.Item(139, 11).Value = 0   # K139
$ws.Cells.Item(139, 12).Value = 0   # L139
$ws.Cells.Item(140, 8).Value = 0   # H140
$ws.Cells.Item(140, 9).Value = 0   # I140
$ws.Cells.Item(140, 10).Value = 0   # J140
$ws.Cells.Item(140, 11).Value = 0   # K140
$ws.Cells.Item(140, 12).Value = 0   # L140
$ws.Cells.Item(141, 8).Value = 0   # H141
$ws.Cells.Item(141, 9).Value = 0   # I141
$ws.Cells.Item(141, 10).Value = 0   # J141
$ws.Cells.Item(141, 11).Value = 0   # K141
$ws.Cells.Item(141, 12).Value = 0   # L141

### Sheet: BSM ###
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117:L117").ClearContents()
$ws.Range("H118:L118").ClearContents()
$ws.Range("H119:L119").ClearContents()
$ws.Range("H120:L120").ClearContents()
$ws.Range("H122:L122").ClearContents()
$ws.Range("H123:L123").ClearContents()
$ws.Range("H124:L124").ClearContents()
$ws.Range("H125:L125").ClearContents()
$ws.Range("H126:L126").ClearContents()
$ws.Range("H127:L127").ClearContents()
$ws.Range("H128:L128").ClearContents()
$ws.Range("H129:L129").ClearContents()
$ws.Range("H130:L130").ClearContents()
$ws.Range("H131:L131").ClearContents()
$ws.Range("H132:L132").ClearContents()
$ws.Range("H133:L133").ClearContents()
$ws.Range("H135:L135").ClearContents()
$ws.Range("H137:L137").ClearContents()
$ws.Range("H138:L138").ClearContents()
$ws.Range("H139:L139").ClearContents()
$ws.Range("H141:L141").ClearContents()
$ws.Range("H134:N134").ClearContents()
$ws.Range("H140:N140").ClearContents()
